# Re-tuned gains to correct formula -.1 -.001 -3;
# added throttle input as an input variable in main()
#
# Appends the new trial-log rows (29-37) to the "Gains_trials" sheet,
# mirroring the existing kp/ki/kd + Notes log layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: note-only row marking the throttle bump to 0.3
$ws.Range("A29").Value = "up throttle to 0.3"

# Row 30: gains tried at that throttle
$ws.Range("A30").Value = -0.1
$ws.Range("B30").Value = 0.001
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = "passed 3 laps ! "

# Row 31: note-only row marking the throttle bump to 0.35
$ws.Range("A31").Value = "up throttle to 0.35"

# Row 32: gains tried at that throttle
$ws.Range("A32").Value = -0.1
$ws.Range("B32").Value = 0.001
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = "passed 3 laps ! "

# Row 33: repeated note confirming throttle still at 0.35
$ws.Range("A33").Value = "up throttle to 0.35"

# Row 34: gains tried again
$ws.Range("A34").Value = -0.1
$ws.Range("B34").Value = 0.001
$ws.Range("C34").Value = 3
$ws.Range("D34").Value = "passed 3 laps ! "

# Row 35: note-only row - bug found, throttle reset to 0.25
$ws.Range("A35").Value = "Found a mistake in my code which, reset throttle to 0.25"

# Row 36: gains tried after the fix, oscillations observed
$ws.Range("A36").Value = -0.1
$ws.Range("B36").Value = 0.001
$ws.Range("C36").Value = 3
$ws.Range("D36").Value = "Crazy Osclilations"

# Row 37: further gains tried, still oscillating
$ws.Range("A37").Value = -0.1
$ws.Range("B37").Value = -0.001
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = "Crazy Osclilations"

# Update workbook view to match the new extent of the log (best effort -
# the runtime may not persist window scroll position on save).
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C38").Select()
